$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old rows 4-13 entirely (not just their contents) so the
# used range / dimension shrinks back down to A1:B3
$ws.Range("A4:B13").EntireRow.Delete()

# New, smaller data set
$ws.Range("A1").Value = 1
$ws.Range("B1").Value = 5

$ws.Range("A2").Value = 3
$ws.Range("B2").Value = 9

$ws.Range("A3").Value = 5
$ws.Range("B3").Value = 25

# Restore the (odd but faithfully reproduced) selection state
$ws.Range("A13").Select()
